$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row for "RM 232" (originally row 26) and "SC 92" (originally row 28 before first deletion, i.e. row 27 after).
# Delete higher row index first so the other row index remains valid.
$ws.Rows.Item(28).Delete()   # SC 92
$ws.Rows.Item(26).Delete()   # RM 232

# Apply per-cell value changes to remaining rows (now renumbered 2..33)
# Row 2: RM 2
$ws.Range("D2").Value = $null
$ws.Range("F2").Value = $null

# Row 4: RM 9
$ws.Range("E4").Value = $null

# Row 5: RM 14
$ws.Range("D5").Value = -14.4

# Row 6: RM 21
$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7

# Row 8: RM 38
$ws.Range("C8").Value = $null
$ws.Range("F8").Value = 17.05

# Row 9: RM 42
$ws.Range("F9").Value = 17.26

# Row 10: RM 52 a
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = -6.1
$ws.Range("F10").Value = 16.43

# Row 11: RM 58
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = 17.65

# Row 12: RM 81
$ws.Range("C12").Value = 12.5
$ws.Range("E12").Value = $null

# Row 13: RM 88
$ws.Range("D13").Value = $null

# Row 14: RM 90
$ws.Range("C14").Value = $null
$ws.Range("E14").Value = -5.4

# Row 16: RM 103
$ws.Range("E16").Value = $null

# Row 17: RM 116
$ws.Range("C17").Value = 11.2
$ws.Range("E17").Value = $null
$ws.Range("F17").Value = $null

# Row 18: RM 120
$ws.Range("C18").Value = 11.5
$ws.Range("F18").Value = $null

# Row 19: RM 125
$ws.Range("C19").Value = $null
$ws.Range("E19").Value = -6.5
$ws.Range("F19").Value = $null

# Row 20: RM 134
$ws.Range("C20").Value = $null

# Row 21: RM 135
$ws.Range("E21").Value = -8.699999999999999

# Row 22: RM 138
$ws.Range("E22").Value = -6.1
$ws.Range("F22").Value = 16.81

# Row 23: RM 140
$ws.Range("C23").Value = 12.2
$ws.Range("F23").Value = $null

# Row 24: RM 142a
$ws.Range("D24").Value = -13.9
$ws.Range("F24").Value = $null

# Row 25: RM 145
$ws.Range("E25").Value = $null

# Row 26: SC 5
$ws.Range("E26").Value = $null

# Row 27: SC 101
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = $null
$ws.Range("E27").Value = $null

# Row 28: SC 105
$ws.Range("D28").Value = $null
$ws.Range("E28").Value = $null
$ws.Range("F28").Value = 17.44

# Row 29: SC 119
$ws.Range("B29").Value = $null
$ws.Range("F29").Value = 18.06

# Row 30: SC 120
$ws.Range("D30").Value = -13.6

# Row 31: SC 132
$ws.Range("E31").Value = -8.1

# Row 32: SC 193
$ws.Range("B32").Value = $null

